$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: "break_on_off" header + 72 data rows (rows 2-73), all 0
# except the rows that mark the end of a training block (19, 37, 54) which get 1.
$ws.Range("L1").Value = "break_on_off"

$breakRows = @(19, 37, 54)
for ($r = 2; $r -le 73; $r++) {
    if ($breakRows -contains $r) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

# Update the view to match: no frozen/scrolled top-left cell, selection
# on the whole new column.
$ws.Range("L1:L73").Select()
